# Applies the two classes of change described by the diff:
#  1. Re-cache the "datetimeFigureOut" date field text (4/17/2022 -> 4/18/2022)
#     on the slide master and every slide layout's Date placeholder.
#  2. Slide 1 subtitle: "(Part 1)" -> "(Section 1)" and merge the trailing
#     "CellRanger (Alignment and Quantification)" runs into a single run.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholder text on the slide master + all custom (slide) layouts.
# ---------------------------------------------------------------------------
$newDate = "4/18/2022"

$containers = @()
$containers += $p.SlideMaster
for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $containers += $p.SlideMaster.CustomLayouts.Item($i)
}

foreach ($container in $containers) {
    $shapes = $container.Shapes
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.Type -eq 14 -and $shp.PlaceholderFormat.Type -eq 16) {
            # ppPlaceholderDate
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 1 ("Subtitle 2" placeholder) text tweaks.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$subtitleShape = $slide1.Shapes.Item(2)
$subtitleRange = $subtitleShape.TextFrame.TextRange

# (Part 1) -> (Section 1)
$currentText = $subtitleRange.Text
$partIdx = $currentText.IndexOf("(Part 1)")
if ($partIdx -ge 0) {
    $partRange = $subtitleRange.Characters($partIdx + 1, 8)
    $partRange.Text = "(Section 1)"
}

# " " + "(Alignment and " + "Quantification)" -> " (Alignment and Quantification)"
$currentText2 = $subtitleRange.Text
$alignIdx = $currentText2.IndexOf(" (Alignment and Quantification)")
if ($alignIdx -ge 0) {
    $alignRange = $subtitleRange.Characters($alignIdx + 1, 32)
    $alignRange.Text = " (Alignment and Quantification)"
}
